$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Header row (row 1): the existing B1/C1/D1 cells were mistakenly
# populated with data-looking values instead of real column labels;
# fix them and extend the header with the new metadata columns E:K
# (matching the pattern used on every other sheet in this workbook). ---
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# Copy the bold/bordered header formatting already used on B1 onto the
# newly added header cells.
$ws.Range("B1").Copy()
$ws.Range("E1:K1").PasteSpecial(-4122)

# --- Data rows 2-8: fill in the new metadata columns E:K that mirror
# the other sheets (property_category, category, date, legislator_name,
# legislator_id, source_file, index). ---
$indexes = @(110, 111, 112, 113, 114, 115, 116)

for ($i = 0; $i -lt 7; $i++) {
    $r = $i + 2
    $ws.Range("E$r").Value = "insurance"
    $ws.Range("F$r").Value = "normal"
    # Force the date to stay a literal text string instead of being
    # auto-converted into a date serial number.
    $ws.Range("G$r").NumberFormat = "@"
    $ws.Range("G$r").Value = "2013-12-24"
    $ws.Range("G$r").NumberFormat = "General"
    $ws.Range("H$r").Value = "段宜康"
    $ws.Range("I$r").Value = 917
    $ws.Range("J$r").Value = "tmpac2a1"
    $ws.Range("K$r").Value = $indexes[$i]
}

# Copy the plain data-row formatting already used on B2 onto the newly
# added data cells.
$ws.Range("B2").Copy()
$ws.Range("E2:K8").PasteSpecial(-4122)
